# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/border/centered
# style used by the rest of row 1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-56): team record for every player row ---
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 99
    $ws.Cells.Item($r, 31).Value = 63
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Done adding Wins/Losses/Ties columns."
